# ---------------------------------------------------------------------------
# Lancers scrape refresh: 2025-11-04 18:25 JST
#
# The scraper re-ran and produced a new snapshot of job listings:
#   - 4 brand-new high-priority AI/Next.js listings at the top
#   - 4 brand-new SRE/OR listings inserted mid-list
#   - all previously-seen listings carried over with a refreshed timestamp
#   - column H (skill summary) widened to fit the new tag text
#
# Simplest faithful way to reproduce this with the COM object model is to
# rebuild the sheet's data region from scratch (clear, then rewrite row by
# row + re-add hyperlinks) rather than trying to replay individual inserts.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop every existing hyperlink (and its relationship part) up front -- the
# URLs/row order are all changing, so there is nothing worth preserving.
$ws.Hyperlinks.Delete()

# Wipe the old data/header region; it is about to be rewritten in full.
$ws.Cells.Clear()

$headers = @('取得日時', 'タイトル', 'カテゴリ', '価格', '締切', 'URL', '優先度スコア', 'スキル概要')
for ($c = 1; $c -le $headers.Count; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Each element: 取得日時, タイトル, カテゴリ, 価格, 締切, URL, 優先度スコア, スキル概要
# スキル概要 is $null where the listing has no skill tags (column H left blank).
$rows = @(
    @('2025-11-04 18:25:38', '医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5416301', 385, '🔥AI,Ai ◆開発 ◇アプリ'),
    @('2025-11-04 18:25:38', 'Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5416328', 310, '🔥AI,Ai'),
    @('2025-11-04 18:25:38', '詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5427010', 245, '🔥Next.js ◆開発,Node.js ◇アプリ'),
    @('2025-11-04 18:25:38', '<Next.js、バックエンド開発> ガントチャートアプリの改修製造', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5427011', 225, '🔥Next.js ◆開発 ◇アプリ'),
    @('2025-11-04 18:25:38', '【急募】システム一元化のための開発', 'システム開発', '5,000,000 円 ~ / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426883', 90, '◆開発'),
    @('2025-11-04 18:25:38', '【急募】WordPressで施設検索サイトのMVPを相談しながら構築いただける方を探しています!', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426483', 65, '◇サイト ○WordPress'),
    @('2025-11-04 18:25:38', 'PHP業務アプリケーションの改修対応', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426598', 58, '◇アプリ ○PHP'),
    @('2025-11-04 18:25:38', '社外エンジニア(WEBサイトやシステムのメンテナンス等の保守/改修等)の募集', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426251', 53, '◇サイト'),
    @('2025-11-04 18:25:38', '【急募】UnityプログラムをiPhoneアプリにコンパイルできる方', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426922', 30, '◇アプリ'),
    @('2025-11-04 18:25:38', '【急募】警備スタッフと各作業現場のマッチングシステム構築依頼', 'システム開発', '1,000,000 円 ~ 3,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426527', 40, $null),
    @('2025-11-04 18:25:38', '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426900', 25, $null),
    @('2025-11-04 18:25:38', '【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426902', 25, $null),
    @('2025-11-04 18:25:38', 'OR(operations research)にて最適化の仕組みの構築(社内常駐)', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5427009', 25, $null),
    @('2025-11-04 18:25:38', 'OR(operations research)にて最適化の仕組みの構築 (リモート)', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5427007', 25, $null),
    @('2025-11-04 18:25:38', '〖リモート可〗Delphiエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5341051', 25, $null),
    @('2025-11-04 18:25:38', 'Amazonセラセンで販売している商品の購入者に対してメッセージを一括送信できるGoogle拡張機能', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426687', 18, $null),
    @('2025-11-04 18:25:38', 'PowerAutomate GoogleドライブからCSVをダウンロードしてヤマトWEBで印刷処理', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5426627', 13, $null)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]   # 取得日時
    $ws.Cells.Item($r, 2).Value = $row[1]   # タイトル
    $ws.Cells.Item($r, 3).Value = $row[2]   # カテゴリ
    $ws.Cells.Item($r, 4).Value = $row[3]   # 価格
    $ws.Cells.Item($r, 5).Value = $row[4]   # 締切
    $ws.Cells.Item($r, 6).Value = $row[5]   # URL (text)
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5]) | Out-Null
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
    $ws.Cells.Item($r, 7).Value = $row[6]   # 優先度スコア (numeric)
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]   # スキル概要 (only when present)
    }
    $r++
}

# Column H widened (17 -> 27 chars) to fit the longer tag strings.
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668
